$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = '2번째 Data'

# Delete old rows 12-21 (articles no longer in the feed)
$ws.Rows("12:21").Delete()

# Overwrite rows 2-11 with the refreshed article feed (newest crawl)
$ws.Range("A2").Value2 = 'China-linked cyber espionage group Velvet Ant exploits Cisco NX-OS zero-day vulnerability to attack… Caution'
$ws.Range("B2").Value2 = 'https://www.dailysecu.com/news/articleView.html?idxno=157415'
$ws.Range("C2").Value2 = 'Velvet Ant, a Chinese-linked cyber espionage group, is exploiting a zero-day vulnerability in Cisco NX-OS software, according to a recent report from cybersecurity firm Sygnia.CVE-2024-20399 The vulnerability being tracked is a command injection vulnerability that could allow an authenticated, local attacker to execute arbitrary commands with root privileges on the underlying operating system of an affected device. This vulnerability, with a CVSS score of 6.0, affects certain vulnerabilities in Cisco NX-OS software. Due to insufficient validation of arguments passed to the configuration CLI command.'
$ws.Range("D2").Value2 = 'Reporter Gil Min-kwon'
$ws.Range("E2").Value2 = '2024-07-03 12:02:26'

$ws.Range("A3").Value2 = 'Kaspersky “Cyber ​​Threat Intelligence (TI) product line and training will be sold as before in the U.S.”'
$ws.Range("B3").Value2 = 'https://www.dailysecu.com/news/articleView.html?idxno=157408'
$ws.Range("C3").Value2 = 'The U.S. Department of Commerce''s Bureau of Industry and Security (BIS) last week announced measures to ban Kaspersky''s U.S. subsidiary from directly or indirectly providing security software in the United States. However, sales and training of Threat Intelligence (TI) products are possible as before. Kaspersky will be prohibited from selling software, excluding TI products and training, to U.S. consumers and companies starting July 20. However, we will continue to provide software and anti-virus signature updates to existing customers until September 29. Additionally, we will provide current individual and business customers with an appropriate replacement software within 100 days to ensure there are no gaps in security protection.'
$ws.Range("D3").Value2 = 'Reporter Gil Min-kwon'
$ws.Range("E3").Value2 = '2024-07-03 09:40:42'

$ws.Range("A4").Value2 = 'Prudential Financial, personal information of 2.5 million people stolen... Alfab ransomware group'
$ws.Range("B4").Value2 = 'https://www.dailysecu.com/news/articleView.html?idxno=157434'
$ws.Range("C4").Value2 = 'Prudential Financial, an American multinational financial services company, was hacked in February and the personal information of more than 2.5 million people was leaked. This figure far exceeds the initially reported number of 36,000. This data breach occurred on February 4, 2024. Prudential Financial detected this the next day and immediately took action. According to Form 8-K filed with the U.S. Securities and Exchange Commission (SEC), the attackers gained access to the company''s administrative and user data, as well as the accounts of employees and contractors. Meanwhile, the BlackCat ransomware group took advantage of this attack.'
$ws.Range("D4").Value2 = 'Reporter Gil Min-kwon'
$ws.Range("E4").Value2 = '2024-07-03 17:15:40'

$ws.Range("A5").Value2 = 'Google conducts kvmCTF to strengthen KVM hypervisor security... Pays $250,000 for zero-day vulnerabilities'
$ws.Range("B5").Value2 = 'https://www.dailysecu.com/news/articleView.html?idxno=157414'
$ws.Range("C5").Value2 = 'Google announced the introduction of a new Vulnerability Reward Program (VRP), kvmCTF, to enhance the security of its Kernel-based Virtual Machine (KVM) hypervisor. First announced in October 2023, the program identifies vulnerabilities within KVM and Aim to correct. KVM plays an important role in consumer and enterprise environments and is used on a variety of platforms, including Android and Google Cloud Platform. KVM is an open source hypervisor with over 17 years of development history that enables Linux to operate as a hypervisor. Google is the master of KVM'
$ws.Range("D5").Value2 = 'Reporter Gil Min-kwon'
$ws.Range("E5").Value2 = '2024-07-03 11:21:41'

$ws.Range("A6").Value2 = 'Soaring ransomware payouts and the limitations of legacy MFA'
$ws.Range("B6").Value2 = 'https://www.dailysecu.com/news/articleView.html?idxno=157423'
$ws.Range("C6").Value2 = 'Cyber ​​security threats have been increasing rapidly recently. According to the Sophos "State of Ransomware 2024" report, average ransom payouts increased 500% last year, reaching an average of $2 million, up from $400,000 in 2023. Additionally, according to the RISK & INSURANCE report, the median ransom demand in 2023 surged from $1.4 million in 2022 to $20 million, and the actual payment amount increased from $335,000 to $6.5 million. ◆Increasing ransomware payouts Factor 1. Targeting of CybercriminalsCybercriminals use maximum'
$ws.Range("D6").Value2 = 'Reporter Gil Min-kwon'
$ws.Range("E6").Value2 = '2024-07-03 15:51:47'

$ws.Range("A7").Value2 = 'Iroun & Company, a generative AI security company, selected as Deeptech Tips Fast Track'
$ws.Range("B7").Value2 = 'https://www.dailysecu.com/news/articleView.html?idxno=157412'
$ws.Range("C7").Value2 = 'Lee Roun & Company (CEO Yoon Du-sik), Korea''s first generative AI security company, was finally selected for the fast track ''Deep Tech Tips'' private investment-led technology startup support program of the Ministry of SMEs and Startups (hereinafter referred to as the Ministry of SMEs and Startups). Lee Roun & Company received seed investment from Mark & ​​Company and the Yoonmin Creative Investment Foundation in February 2024, and was selected as DeepTech TIPS after being recommended by TIPS operator Mark & ​​Company. DeepTech TIPS is a promising start in the top 10 new industries. As part of the Super Gap Startup 1000+ project announced by the Ministry of SMEs and Startups to select and foster businesses, selected companies will receive KRW 1.5 billion in research and development (R&D) expenses, as well as commercialization and overseas support.'
$ws.Range("D7").Value2 = 'Reporter Gil Min-kwon'
$ws.Range("E7").Value2 = '2024-07-03 10:48:24'

$ws.Range("A8").Value2 = 'JSI Lab launches X-ITM, Auto ML-based threat response solution'
$ws.Range("B8").Value2 = 'https://www.dailysecu.com/news/articleView.html?idxno=157417'
$ws.Range("C8").Value2 = 'JSI LAB (CEO Taeyong Ha), a company specializing in security AI, announced that it has launched X-ITM 1.0, an Auto ML-based threat response solution. JSI LAB conducts Auto ML and AI work. In the security field, predicting security events based on transformers, recognizing hacking through deep learning-based image recognition, predicting security incidents through statistical-based predictions, predicting unknown or unrecognized patterns through unsupervised learning, and predicting patterns that are unknown or unrecognized through existing ESM or Load patterns into SIEM and improve the accuracy of security events through supervised learning using existing patterns and labeling.'
$ws.Range("D8").Value2 = 'Reporter Gil Min-kwon'
$ws.Range("E8").Value2 = '2024-07-03 14:00:04'

$ws.Range("A9").Value2 = 'RaonSecure wins contract for ‘mobile ID construction and expansion service’ project'
$ws.Range("B9").Value2 = 'https://www.dailysecu.com/news/articleView.html?idxno=157418'
$ws.Range("C9").Value2 = 'Raon Secure (CEO Soon-Hyung Lee, Jeong-Ah Lee) announced that it has won an order for the ‘Mobile ID Construction and Expansion Service’ project ordered by the Korea Minting and Security Minting Corporation to build a blockchain-based mobile resident registration card issuance system by the Ministry of the Interior and Safety. Raon Secure is issuing this mobile resident registration card. We provide OmniOne Enterprise, a blockchain-based digital ID platform, for system development and are responsible for the mobile security sector required when issuing resident registration cards in a mobile environment. RaonSecure''s blockchain-based digital technology has previously been applied to mobile driver''s licenses and mobile national veterans registration cards. OmniOne Enterprise, an identity platform, is a FIDO (biometric'
$ws.Range("D9").Value2 = 'Reporter Gil Min-kwon'
$ws.Range("E9").Value2 = '2024-07-03 14:16:46'

$ws.Range("A10").Value2 = 'AhnLab selected as a supplier for the ‘2024 ICT Small and Medium Business Information Protection Support Project’...Recruiting demand companies'
$ws.Range("B10").Value2 = 'https://www.dailysecu.com/news/articleView.html?idxno=157420'
$ws.Range("C10").Value2 = 'AhnLab (CEO Seok-gyun Kang) announced that it has been selected as a security solution supplier for the ''2024 ICT Small and Medium Business Information Protection Support Project (hereinafter referred to as ''Support Project'') and is recruiting companies in demand. This support project includes ◇ information security consulting and security solution support ◇ cloud It is divided into two areas, including support for basic security services (SECaaS). AhnLab provides SaaS (Software as a Service) type security management solution for small and medium-sized businesses across each field, including the ''AhnLab Office Security'' product line (3 types); Network security solution ‘AhnLab TrusGuard’ 50B-70B model'
$ws.Range("D10").Value2 = 'Reporter Gil Min-kwon'
$ws.Range("E10").Value2 = '2024-07-03 14:40:29'

$ws.Range("A11").Value2 = 'KISA-Korea Aerospace Research Institute signs business agreement to strengthen aerospace security'
$ws.Range("B11").Value2 = 'https://www.dailysecu.com/news/articleView.html?idxno=157428'
$ws.Range("C11").Value2 = 'Korea Internet & Security Agency (President Lee Sang-jung, hereinafter referred to as KISA) signed a business agreement (MOU) to strengthen security with the Korea Aerospace Research Institute (President Lee Sang-ryul, hereinafter referred to as ''Aerospace Research Institute'') at the Korea Aerospace Research Institute''s Daejeon headquarters on Wednesday, July 3. KISA made efforts to spread security awareness by collaborating with private aerospace industry officials, such as holding an ''Aerospace and Aviation Convergence Security Seminar'' in October last year to share the latest aerospace security trends and policy directions. This agreement Accordingly, both organizations will strengthen the security of aerospace companies and products, develop and verify aerospace security models, train aerospace security personnel and raise awareness, and'
$ws.Range("D11").Value2 = 'Reporter Gil Min-kwon'
$ws.Range("E11").Value2 = '2024-07-03 16:47:35'

